{"js": "// Update the worksheet date and every two-digit-by-two-digit multiplication\n// problem/answer pair in the table, per the target revision.\nconst replacements = [\n  [\"2026-02-02 Monday\", \"2026-02-03 Tuesday\"],\n  [\"36\u00d753=1908\", \"21\u00d775=1575\"],\n  [\"73\u00d751=3723\", \"14\u00d727=378\"],\n  [\"73\u00d793=6789\", \"98\u00d734=3332\"],\n  [\"27\u00d781=2187\", \"38\u00d716=608\"],\n  [\"41\u00d721=861\", \"79\u00d786=6794\"],\n  [\"87\u00d779=6873\", \"58\u00d742=2436\"],\n  [\"26\u00d780=2080\", \"64\u00d791=5824\"],\n  [\"23\u00d775=1725\", \"72\u00d744=3168\"],\n  [\"97\u00d794=9118\", \"75\u00d714=1050\"],\n  [\"27\u00d747=1269\", \"68\u00d724=1632\"],\n  [\"97\u00d711=1067\", \"57\u00d719=1083\"],\n  [\"31\u00d788=2728\", \"36\u00d749=1764\"],\n  [\"66\u00d786=5676\", \"67\u00d723=1541\"],\n  [\"21\u00d725=525\", \"43\u00d758=2494\"],\n  [\"46\u00d736=1656\", \"95\u00d765=6175\"],\n  [\"92\u00d712=1104\", \"22\u00d755=1210\"],\n  [\"43\u00d712=516\", \"69\u00d782=5658\"],\n  [\"75\u00d744=3300\", \"48\u00d721=1008\"],\n  [\"38\u00d715=570\", \"41\u00d733=1353\"],\n  [\"39\u00d743=1677\", \"76\u00d737=2812\"],\n  [\"91\u00d756=5096\", \"36\u00d750=1800\"],\n  [\"69\u00d726=1794\", \"44\u00d772=3168\"],\n  [\"43\u00d766=2838\", \"34\u00d753=1802\"],\n  [\"45\u00d756=2520\", \"69\u00d711=759\"],\n  [\"93\u00d793=8649\", \"85\u00d760=5100\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  for (const item of results.items) {\n    item.insertText(newText, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "# Update the worksheet date and every two-digit-by-two-digit multiplication\n# problem/answer pair in the table, per the target revision.\n$d = $word.ActiveDocument\n\n$pairs = @(\n    @(\"2026-02-02 Monday\", \"2026-02-03 Tuesday\"),\n    @(\"36\u00d753=1908\", \"21\u00d775=1575\"),\n    @(\"73\u00d751=3723\", \"14\u00d727=378\"),\n    @(\"73\u00d793=6789\", \"98\u00d734=3332\"),\n    @(\"27\u00d781=2187\", \"38\u00d716=608\"),\n    @(\"41\u00d721=861\", \"79\u00d786=6794\"),\n    @(\"87\u00d779=6873\", \"58\u00d742=2436\"),\n    @(\"26\u00d780=2080\", \"64\u00d791=5824\"),\n    @(\"23\u00d775=1725\", \"72\u00d744=3168\"),\n    @(\"97\u00d794=9118\", \"75\u00d714=1050\"),\n    @(\"27\u00d747=1269\", \"68\u00d724=1632\"),\n    @(\"97\u00d711=1067\", \"57\u00d719=1083\"),\n    @(\"31\u00d788=2728\", \"36\u00d749=1764\"),\n    @(\"66\u00d786=5676\", \"67\u00d723=1541\"),\n    @(\"21\u00d725=525\", \"43\u00d758=2494\"),\n    @(\"46\u00d736=1656\", \"95\u00d765=6175\"),\n    @(\"92\u00d712=1104\", \"22\u00d755=1210\"),\n    @(\"43\u00d712=516\", \"69\u00d782=5658\"),\n    @(\"75\u00d744=3300\", \"48\u00d721=1008\"),\n    @(\"38\u00d715=570\", \"41\u00d733=1353\"),\n    @(\"39\u00d743=1677\", \"76\u00d737=2812\"),\n    @(\"91\u00d756=5096\", \"36\u00d750=1800\"),\n    @(\"69\u00d726=1794\", \"44\u00d772=3168\"),\n    @(\"43\u00d766=2838\", \"34\u00d753=1802\"),\n    @(\"45\u00d756=2520\", \"69\u00d711=759\"),\n    @(\"93\u00d793=8649\", \"85\u00d760=5100\")\n)\n\nforeach ($pair in $pairs) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n\n    $range = $d.Content\n    $range.Find.ClearFormatting()\n    $range.Find.Execute($oldText, $false, $false, $false, $false, $false, $true, 1, $false, $newText, 2) | Out-Null\n}\n"}
